$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FilePath -> CloneScene dir, ID -> 3, SceneName -> clone
$ws.Range("A2").Value = "../../NFDataCfg/Ini/NFZoneServer/Scene/CloneScene/"
$ws.Range("B2").Value = "3"
$ws.Range("F2").Value = "clone"

# Row 3: SceneName -> newscene
$ws.Range("F3").Value = "newscene"

# Row 4: SceneName -> newscene
$ws.Range("F4").Value = "newscene"

# Update the active cell selection to match the author's last position
$ws.Range("H8").Select()
